$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the cached "datetimeFigureOut" field text from 4/25/2019 to
#    4/28/2019 everywhere it appears: the notes master, every slide
#    layout, and the slide master itself.
# ---------------------------------------------------------------------
$oldDate = "4/25/2019"
$newDate = "4/28/2019"

# Notes master
$nm = $p.NotesMaster
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $sh = $nm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
$sm = $p.SlideMaster
for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Every slide layout
for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $cl = $sm.CustomLayouts.Item($li)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $sh = $cl.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 38 ("Linear Regression (OLS) Assumptions"): tidy up the
#    '‘Omitted Variable bias’ ' phrase to '‘Omitted Variable bias '
#    (drop the closing curly quote and reflow the bold run boundary).
# ---------------------------------------------------------------------
$s38 = $p.Slides.Item(38)
$shp = $s38.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$txt = $tr.Text
$markerIdx = $txt.IndexOf("``Omitted Variable bias' ")
if ($markerIdx -ge 0) {
    $sub = $tr.Characters($markerIdx + 1 + 9, 15)
    $sub.Text = "Variable bias "
}
